$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 536, pushing the existing rows (536-563) down to (537-564).
$ws.Rows.Item(536).Insert()

# Populate the newly inserted row 536 with a new weekly price observation.
# Most attributes mirror the row that used to occupy 536 (now at 537); only
# the date (D) and volume (J) differ for this new record.
$ws.Range("A536").Value = 8
$ws.Range("B536").Value = "Terminal La Palmera de La Serena"
$ws.Range("C536").Value = "Coquimbo"
$ws.Range("D536").Value = 45147
$ws.Range("E536").Value = 4
$ws.Range("F536").Value = 100112032
$ws.Range("G536").Value = "Zapallo italiano"
$ws.Range("H536").Value = "Sin especificar"
$ws.Range("I536").Value = "Primera"
$ws.Range("J536").Value = 500
$ws.Range("K536").Value = 14000
$ws.Range("L536").Value = 15000
$ws.Range("M536").Value = 14500
$ws.Range("N536").Value = "$/caja 50 unidades"
$ws.Range("O536").Value = "Región de Arica y Parinacota"
$ws.Range("P536").Value = 290
$ws.Range("Q536").Value = 50
$ws.Range("R536").Value = "Hortaliza"

# Keep the date formatting consistent with the rest of column D.
$ws.Range("D536").NumberFormat = $ws.Range("D537").NumberFormat
